$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H holds the "Absent" count for each date row, derived from
# column D ("Total Attendance Count"): absent (1) when there was no
# attendance recorded that day (D = 0), otherwise present (0).
# This consolidates the report by filling in the previously
# inconsistent/missing values in column H.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    if ($total -eq $null) { $total = 0 }

    if ($total -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
